# Generate Report for Handoff
# - Priority for the "Ready for handoff" rows (4-7) moves from "low" to "ht"
#   on both the zh-cn and de-de localization-status sheets.
# - Latest Handoff Datetime for those same rows is refreshed to reflect the
#   newly (re)generated handoff xliff files.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4").Value = "2016-08-28 00:29:29"
$zhcn.Range("H5").Value = "2016-08-28 00:29:29"
$zhcn.Range("H6").Value = "2016-08-28 00:29:29"
$zhcn.Range("H7").Value = "2016-08-28 00:29:29"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4").Value = "2016-08-28 00:29:34"
$dede.Range("H5").Value = "2016-08-28 00:29:34"
$dede.Range("H6").Value = "2016-08-28 00:29:34"
$dede.Range("H7").Value = "2016-08-28 00:29:34"
